$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Jalen Green","PG,SG","Houston Rockets"),
    @("Jalen Suggs","PG,SG","Orlando Magic"),
    @("Clint Capela","C","Atlanta Hawks"),
    @("Jakob Poeltl","C","Toronto Raptors"),
    @("Rudy Gobert","C","Minnesota Timberwolves"),
    @("Naz Reid","PF,C","Minnesota Timberwolves"),
    @("Jerami Grant","SF,PF","Portland Trail Blazers"),
    @("Deni Avdija","SF,PF","Portland Trail Blazers"),
    @("Pascal Siakam","SF,PF","Indiana Pacers"),
    @("Chris Paul","PG","San Antonio Spurs"),
    @("Dejounte Murray","PG,SG","New Orleans Pelicans"),
    @("Russell Westbrook","PG","Denver Nuggets"),
    @("Nikola Jokic","C","Denver Nuggets"),
    @("Jaylen Brown","SG,SF","Boston Celtics"),
    @("Ayo Dosunmu","SG,SF","Chicago Bulls"),
    @("Paolo Banchero","SF,PF","Orlando Magic"),
    @("Chet Holmgren","PF,C","Oklahoma City Thunder"),
    @("Bogdan Bogdanovic","SG,SF","Atlanta Hawks")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
